# Auto-applied market-data refresh for Leve profit tracker sheets
# (values sourced from the latest Universalis price snapshot)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 1325
$ws.Cells.Item(2, 10).Value = 650
$ws.Cells.Item(2, 12).Value = 650
$ws.Cells.Item(2, 14).Value = -876
$ws.Cells.Item(58, 8).Value = 1282.0952
$ws.Cells.Item(58, 9).Value = 769.05884
$ws.Cells.Item(58, 10).Value = 3462.5
$ws.Cells.Item(58, 11).Value = 2307.17652
$ws.Cells.Item(58, 12).Value = 10387.5
$ws.Cells.Item(58, 13).Value = -2157.17652
$ws.Cells.Item(58, 14).Value = -10687.5
$ws.Cells.Item(88, 8).Value = 2133.8333
$ws.Cells.Item(88, 9).Value = 3000
$ws.Cells.Item(88, 10).Value = 1960.6
$ws.Cells.Item(88, 11).Value = 3000
$ws.Cells.Item(88, 12).Value = 1960.6
$ws.Cells.Item(88, 13).Value = -2594
$ws.Cells.Item(88, 14).Value = -2772.6
$ws.Cells.Item(91, 8).Value = 2133.8333
$ws.Cells.Item(91, 9).Value = 3000
$ws.Cells.Item(91, 10).Value = 1960.6
$ws.Cells.Item(91, 11).Value = 3000
$ws.Cells.Item(91, 12).Value = 1960.6
$ws.Cells.Item(91, 13).Value = -1596
$ws.Cells.Item(91, 14).Value = -4768.6
$ws.Cells.Item(110, 8).Value = 50702
$ws.Cells.Item(110, 10).Value = 50702
$ws.Cells.Item(110, 12).Value = 50702
$ws.Cells.Item(110, 14).Value = -58882
$ws.Cells.Item(111, 8).Value = 3111
$ws.Cells.Item(111, 9).Value = 1000
$ws.Cells.Item(111, 10).Value = 3638.75
$ws.Cells.Item(111, 11).Value = 3000
$ws.Cells.Item(111, 12).Value = 10916.25
$ws.Cells.Item(111, 13).Value = 67
$ws.Cells.Item(111, 14).Value = -17050.25
$ws.Cells.Item(112, 8).Value = 2987.44
$ws.Cells.Item(112, 9).Value = 1250
$ws.Cells.Item(112, 10).Value = 3318.3809
$ws.Cells.Item(112, 11).Value = 3750
$ws.Cells.Item(112, 12).Value = 9955.1427
$ws.Cells.Item(112, 13).Value = -2642
$ws.Cells.Item(112, 14).Value = -12171.1427
$ws.Cells.Item(116, 8).Value = 5834
$ws.Cells.Item(116, 9).Value = 5914.1665
$ws.Cells.Item(116, 10).Value = 5696.5713
$ws.Cells.Item(116, 11).Value = 5914.1665
$ws.Cells.Item(116, 12).Value = 5696.5713
$ws.Cells.Item(116, 13).Value = -2472.1665
$ws.Cells.Item(116, 14).Value = -12580.5713
$ws.Cells.Item(125, 8).Value = 8290.25
$ws.Cells.Item(125, 9).Value = 8290.25
$ws.Cells.Item(125, 11).Value = 74612.25
$ws.Cells.Item(125, 13).Value = -72152.25
$ws.Cells.Item(138, 8).Value = 2892.85
$ws.Cells.Item(138, 10).Value = 2766.5615
$ws.Cells.Item(138, 12).Value = 8299.684499999999
$ws.Cells.Item(138, 14).Value = -18579.6845

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5840.953
$ws.Cells.Item(32, 9).Value = 3677.875
$ws.Cells.Item(32, 11).Value = 3677.875
$ws.Cells.Item(32, 13).Value = -3390.875
$ws.Cells.Item(45, 8).Value = 5997.364
$ws.Cells.Item(45, 9).Value = 5502.5
$ws.Cells.Item(45, 11).Value = 5502.5
$ws.Cells.Item(45, 13).Value = -5125.5
$ws.Cells.Item(110, 8).Value = 5057835.5
$ws.Cells.Item(110, 9).Value = 8270883
$ws.Cells.Item(110, 10).Value = 8761.143
$ws.Cells.Item(110, 11).Value = 8270883
$ws.Cells.Item(110, 12).Value = 8761.143
$ws.Cells.Item(110, 13).Value = -8268838
$ws.Cells.Item(110, 14).Value = -12851.143
$ws.Cells.Item(122, 8).Value = 606369.1
$ws.Cells.Item(122, 9).Value = 928574.5600000001
$ws.Cells.Item(122, 10).Value = 7987.5713
$ws.Cells.Item(122, 11).Value = 2785723.68
$ws.Cells.Item(122, 12).Value = 23962.7139
$ws.Cells.Item(122, 13).Value = -2783273.68
$ws.Cells.Item(122, 14).Value = -28862.7139
$ws.Cells.Item(139, 8).Value = 84999
$ws.Cells.Item(139, 10).Value = 84999
$ws.Cells.Item(139, 12).Value = 84999
$ws.Cells.Item(139, 14).Value = -95279

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 1970.5518
$ws.Cells.Item(105, 9).Value = 1851.7693
$ws.Cells.Item(105, 11).Value = 1851.7693
$ws.Cells.Item(105, 13).Value = -104.7692999999999
$ws.Cells.Item(134, 8).Value = 46493.723
$ws.Cells.Item(134, 9).Value = 48881.566
$ws.Cells.Item(134, 11).Value = 146644.698
$ws.Cells.Item(134, 13).Value = -144109.698

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 1851.9445
$ws.Cells.Item(129, 10).Value = 2627.3333
$ws.Cells.Item(129, 12).Value = 7881.999899999999
$ws.Cells.Item(129, 14).Value = -17881.9999
$ws.Cells.Item(132, 8).Value = 2116.1538
$ws.Cells.Item(132, 10).Value = 2928.4
$ws.Cells.Item(132, 12).Value = 26355.6
$ws.Cells.Item(132, 14).Value = -31415.6
$ws.Cells.Item(139, 8).Value = 9065.385
$ws.Cells.Item(139, 9).Value = 12641.667
$ws.Cells.Item(139, 11).Value = 37925.001
$ws.Cells.Item(139, 13).Value = -32785.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(14, 8).Value = 5000020
$ws.Cells.Item(14, 9).Value = 5000020
$ws.Cells.Item(14, 11).Value = 5000020
$ws.Cells.Item(14, 13).Value = -4999852
$ws.Cells.Item(52, 8).Value = 14436.5
$ws.Cells.Item(52, 10).Value = 14436.5
$ws.Cells.Item(52, 12).Value = 14436.5
$ws.Cells.Item(52, 14).Value = -14954.5
$ws.Cells.Item(70, 8).Value = 22704.883
$ws.Cells.Item(70, 10).Value = 25248.25
$ws.Cells.Item(70, 12).Value = 25248.25
$ws.Cells.Item(70, 14).Value = -25788.25
$ws.Cells.Item(73, 8).Value = 22704.883
$ws.Cells.Item(73, 10).Value = 25248.25
$ws.Cells.Item(73, 12).Value = 25248.25
$ws.Cells.Item(73, 14).Value = -27120.25
$ws.Cells.Item(122, 8).Value = 5092.2163
$ws.Cells.Item(122, 9).Value = 3168.6538
$ws.Cells.Item(122, 10).Value = 9638.817999999999
$ws.Cells.Item(122, 11).Value = 9505.9614
$ws.Cells.Item(122, 12).Value = 28916.454
$ws.Cells.Item(122, 13).Value = -7055.9614
$ws.Cells.Item(122, 14).Value = -33816.454

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1913104.5
$ws.Cells.Item(7, 9).Value = 3974786.2
$ws.Cells.Item(7, 11).Value = 3974786.2
$ws.Cells.Item(7, 13).Value = -3974674.2
$ws.Cells.Item(16, 8).Value = 50002164
$ws.Cells.Item(16, 9).Value = 58825868
$ws.Cells.Item(16, 11).Value = 58825868
$ws.Cells.Item(16, 13).Value = -58825698
$ws.Cells.Item(22, 8).Value = 21742532
$ws.Cells.Item(22, 10).Value = 66672532
$ws.Cells.Item(22, 12).Value = 66672532
$ws.Cells.Item(22, 14).Value = -66673122
$ws.Cells.Item(24, 8).Value = 62111.11
$ws.Cells.Item(24, 9).Value = 20000
$ws.Cells.Item(24, 10).Value = 83166.664
$ws.Cells.Item(24, 11).Value = 20000
$ws.Cells.Item(24, 12).Value = 83166.664
$ws.Cells.Item(24, 13).Value = -19657
$ws.Cells.Item(24, 14).Value = -83852.664
$ws.Cells.Item(27, 8).Value = 21742532
$ws.Cells.Item(27, 10).Value = 66672532
$ws.Cells.Item(27, 12).Value = 66672532
$ws.Cells.Item(27, 14).Value = -66672746
$ws.Cells.Item(40, 8).Value = 1686837.8
$ws.Cells.Item(40, 9).Value = 3749.75
$ws.Cells.Item(40, 10).Value = 5359030
$ws.Cells.Item(40, 11).Value = 3749.75
$ws.Cells.Item(40, 12).Value = 5359030
$ws.Cells.Item(40, 13).Value = -3613.75
$ws.Cells.Item(40, 14).Value = -5359302
$ws.Cells.Item(93, 8).Value = 37044490
$ws.Cells.Item(93, 9).Value = 62503780
$ws.Cells.Item(93, 11).Value = 62503780
$ws.Cells.Item(93, 13).Value = -62502532
$ws.Cells.Item(126, 8).Value = 1913104.5
$ws.Cells.Item(126, 9).Value = 3974786.2
$ws.Cells.Item(126, 11).Value = 11924358.6
$ws.Cells.Item(126, 13).Value = -11921888.6
$ws.Cells.Item(132, 8).Value = 12999.76
$ws.Cells.Item(132, 9).Value = 6282.6665
$ws.Cells.Item(132, 10).Value = 19200.154
$ws.Cells.Item(132, 11).Value = 18847.9995
$ws.Cells.Item(132, 12).Value = 57600.462
$ws.Cells.Item(132, 13).Value = -16317.9995
$ws.Cells.Item(132, 14).Value = -62660.462

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3, 8).Value = 9346668
$ws.Cells.Item(3, 9).Value = 14000000
$ws.Cells.Item(3, 11).Value = 14000000
$ws.Cells.Item(3, 13).Value = -13999886
$ws.Cells.Item(7, 8).Value = 10005
$ws.Cells.Item(7, 10).Value = 10005
$ws.Cells.Item(7, 12).Value = 10005
$ws.Cells.Item(7, 14).Value = -10231
$ws.Cells.Item(19, 8).Value = 15229
$ws.Cells.Item(19, 9).Value = 15000
$ws.Cells.Item(19, 10).Value = 15305.333
$ws.Cells.Item(19, 11).Value = 15000
$ws.Cells.Item(19, 12).Value = 15305.333
$ws.Cells.Item(19, 13).Value = -14826
$ws.Cells.Item(19, 14).Value = -15653.333
$ws.Cells.Item(41, 8).Value = 16619.334
$ws.Cells.Item(41, 10).Value = 15944
$ws.Cells.Item(41, 12).Value = 15944
$ws.Cells.Item(41, 14).Value = -16724
$ws.Cells.Item(45, 8).Value = 13669.4
$ws.Cells.Item(45, 10).Value = 13669.4
$ws.Cells.Item(45, 12).Value = 13669.4
$ws.Cells.Item(45, 14).Value = -14651.4
$ws.Cells.Item(113, 8).Value = 707.6923
$ws.Cells.Item(113, 9).Value = 220.1
$ws.Cells.Item(113, 11).Value = 660.3
$ws.Cells.Item(113, 13).Value = 1509.7
$ws.Cells.Item(141, 8).Value = 74000
$ws.Cells.Item(141, 10).Value = 74000
$ws.Cells.Item(141, 12).Value = 74000
$ws.Cells.Item(141, 14).Value = -84360
